$wb = $excel.ActiveWorkbook

# --- Sheet "Test Cases": clear the stray Result value in F2 ---
$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Select()
$wsTestCases.Cells.Item(2, 6).ClearContents()
$wsTestCases.Range("D2").Select()

# --- Sheet "VerifyStartEndDateValidation": fix the typo'd message, update
#     the demo date/version-lead data, and clear the stale Result/Comments ---
$wsValidation = $wb.Worksheets.Item("VerifyStartEndDateValidation")
$wsValidation.Select()

# Expected message column (I2): fix "Dat" -> "Date" typo
$wsValidation.Cells.Item(2, 9).Value = "End Date should always be greater or equal to the Start Date!"

# End Date column (G2): 24/03/2017 -> 24/04/2017 (leading ' keeps it text, matching existing style)
$wsValidation.Cells.Item(2, 7).Value = "'24/04/2017"

# Version Lead column (H2): Test User -> test demo
$wsValidation.Cells.Item(2, 8).Value = "'test demo"

# Result column (K2) and Comments column (L2): clear stale FAIL / expected-message text
$wsValidation.Cells.Item(2, 11).ClearContents()
$wsValidation.Cells.Item(2, 12).ClearContents()

$wsValidation.Range("D2").Select()
